$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Row 9 (issue #6): status changed from "N/S" to "In progress"
$ws.Range("C9").Value = "In progress"

# 2. Row 12 (issue #9): new defect/fix entry filled in (was blank).
#    First copy the formatting (fill/alignment) from row 7, which uses the
#    same "FIX" row style, then set the row height and the cell values.
$ws.Range("A7:J7").Copy() | Out-Null
$ws.Range("A12:J12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Rows.Item(12).RowHeight = 105

$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "FIX"
$ws.Range("D12").Value = "1.3.3"
$ws.Range("E12").Value = "1.3.3"
$ws.Range("F12").Value = "Any"
$ws.Range("G12").Value = "CONT mode with delay"
$ws.Range("H12").Value = "Same as 1.  Fixing it broke Counin = 0 mode!"
$ws.Range("I12").Value = "seems the two CONT modes do not like each other…  With/without delay = mutually exclusive."
$ws.Range("J12").Value = "Redesign.  Now three Strips Modes to avoiud confusion when using Delay to define which CONT mode we are in.  MAN: manual = stops at the end of each strip;  STEP: Delay at end of each strip;  CONT : does not stop - beeps at 0.5s defore the end of each strip."

# 3. Selection moved to C10 (as last saved by the author)
$ws.Range("C10").Select() | Out-Null
